$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: create row 48 by copying row 47 down (still "Buy / IN PROGRESS" at this point) ---
$null = $ws.Range("A47:L47").Copy($ws.Range("A48:L48"))
$ws.Rows.Item(48).RowHeight = 14.25
# row 47 had no content in J/L -- the ranged copy still materializes blank cells there, drop them
$null = $ws.Range("J48").ClearContents()
$null = $ws.Range("L48").ClearContents()

# --- Step 2: finish off row 47 (the original buy order becomes DONE, gets its fee/finalized date) ---
$ws.Range("J47").Value = "0.31500000 XRP (0.15%)"
$ws.Range("H47").Value = "DONE"
$ws.Range("I47").Value = 42863.315740740742

# --- Step 3: fix up row 48 so it reflects the new sell order ---
$ws.Range("A48").Value = 42863.375740740739

# "Sell", red rich text -- copy from an existing "Sell" cell so the run formatting matches
$null = $ws.Range("B42").Copy($ws.Range("B48"))

# D48 is a numeric-looking piece of text; enter it via a formula (text concatenation keeps
# it a string) then flatten the formula down to a plain cached value.
$ws.Range("D48").Formula = '="           0.18590000"&CHAR(10)'
$null = $ws.Range("D48").Copy()
$null = $ws.Range("D48").PasteSpecial(-4163)

$ws.Range("E48").Value = "         0.197USDT"
$ws.Range("F48").Value = "         209 XRP"

# C48, G48, H48, I48, K48 already came over correctly from the Step 1 row copy

# --- Step 4: update the active selection to match the authored edit ---
$null = $ws.Range("J52").Select()
